# remove crowpasstrail from sde_common
# The "crowpasstrail" entry lives in the row whose Nick (col B) = "crowpasstrail"
# and whose Path (col C) ends with "sde_common.DBO.crowpasstrail".
# Find that row dynamically (rather than hard-coding row 216) and delete it
# entirely, letting Excel re-flow the formulas / autofilter / dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = $ws.Cells(1048576, 2).End(-4162).Row  # xlUp = -4162, search up column B

$targetRow = 0
for ($r = 2; $r -le $lastRow; $r++) {
    $nick = $ws.Cells.Item($r, 2).Value2
    if ($nick -eq "crowpasstrail") {
        $targetRow = $r
        break
    }
}

if ($targetRow -gt 0) {
    $ws.Rows.Item($targetRow).Delete() | Out-Null

    # Mimic the selection state left behind after deleting the row in the UI
    $ws.Cells.Item($targetRow, 1).EntireRow.Select() | Out-Null

    $lastRow = $lastRow - 1
    $lastCol = "F"

    # Re-apply AutoFilter so its range shrinks to match the new data extent
    $ws.AutoFilterMode = $false
    $ws.Range("A1:" + $lastCol + $lastRow).AutoFilter() | Out-Null

    # Keep the hidden _FilterDatabase defined name in sync with the filter range
    foreach ($n in $wb.Names) {
        if ($n.Name -eq "Sheet1!_FilterDatabase") {
            $n.RefersTo = "=Sheet1!`$A`$1:`$" + $lastCol + "`$" + $lastRow
        }
    }
}

$wb.Save()
